$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new row of data (row 4) ---
$ws.Range("B4").Value = 10.24
$ws.Range("C4").Value = 8.57
$ws.Range("D4").Value = 10.63
$ws.Range("E4").Value = 12.45
$ws.Range("F4").Value = 10.050000000000001
$ws.Range("G4").Value = 11.88
$ws.Range("H4").Value = 11.44
$ws.Range("I4").Value = 12.93
$ws.Range("J4").Value = 10.02
$ws.Range("K4").Value = 10.4
$ws.Range("L4").Value = 10.87
$ws.Range("M4").Value = 9.69

# --- Re-apply borders across the whole table (B1:M4), as happens when the
#     user selects the range and presses the "All Borders" button ---

# Row 1 (header row): thin vertical separators between header cells, with a
# medium (heavier) box around the whole header row.
$row1 = $ws.Range("B1:M1")
$row1.Borders.LineStyle = 1
$row1.Borders.Weight = -4138
$row1.Borders.Color = 0
$row1.Borders.Item(11).Weight = 2
$row1.Borders.Item(11).Color = 0

# Row 2: thin grid on left/right/bottom, no line between header and this row's
# top (the header's own medium bottom border already separates them).
$row2 = $ws.Range("B2:M2")
$row2.Borders.LineStyle = 1
$row2.Borders.Weight = 2
$row2.Borders.Color = 0
$row2.Borders.Item(8).LineStyle = -4142

# Row 3: full thin grid box on every cell.
$row3 = $ws.Range("B3:M3")
$row3.Borders.LineStyle = 1
$row3.Borders.Weight = 2
$row3.Borders.Color = 0

# Row 4 (new row): full thin grid box on every cell, matching row 3.
$row4 = $ws.Range("B4:M4")
$row4.Borders.LineStyle = 1
$row4.Borders.Weight = 2
$row4.Borders.Color = 0

# --- Misc view/formatting touch-ups seen in the authored workbook ---
$ws.Rows("1").RowHeight = 16
$ws.Range("B1:M1").Borders.Item(9).LineStyle = 1
$ws.Range("B1:M1").Borders.Item(9).Weight = -4138
$ws.Range("E12").Select
